$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '281474992433979-1743538267069'
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = 'Harsh Brake'
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '2025-04-01T14:11:07.069'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '281474992433979'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '131'
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '51834059'
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = 'DANIEL IÑIGUEZ'
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 20.71519898
$ws.Range("I2").Value = -103.38799362
$ws.Range("J2").Value = 0.7422594428062439
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = 'No video URL'
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = 'No video URL'
$ws.Range("L2").Style = "Normal"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = '281474991205262-1743535592415'
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = 'Forward Collision Warning'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '2025-04-01T13:26:32.415'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '281474991205262'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '132'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '52215867'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = 'EMMANUEL SALCEDO'
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 20.65233998
$ws.Range("I3").Value = -103.3124147
$ws.Range("J3").Value = 0
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEF5NZIN5%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T160057Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG0aCXVzLXdlc3QtMiJHMEUCICQZwBitIu4Neglo3bHNBCT1E8PvFZ%2FQxlTfmkY3HJ4KAiEAg6tKCsnDIROPImOf76nbWN6GO9m4dMJNp1nqHUL9bjsq5gMI1v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDHLIT5JpxqhsIaJq6iq6AzJNAysmsomDmw%2B4G9FYLibD3eb%2Fe704hzNOoDYLc%2B7TRT4zkrhXlK2nwowzVJQdWgN%2BaJw0CXT3sCavmv3mHNKRaGbJ8k8xsskjo05Ge%2FGEc88ZMGDiGcWroBAlx3QvC7B0uSFXY%2BnyC28u6VSxj3PyAkhPd45oMhiaDJqi6b3zmq19EpJfVGdW7F7ihIPserIIo6g9at3uD3uVKeAxlKh%2B0aNHIEeb3QsUpzjUJLT0k3Yr3RzlFwWC%2B8a4kmAFK12kmir3%2FnkP9lPSD8%2FJ3zeT3xqrzaQ%2BhJwxI9IcJ88lWHMpHTNTipcI2XuPlee9kENxebUJehDSPNGKYmho2UxmXF79E0aLRDyC6c8FOpoli1LX7%2ByPQjR8%2BphqtIewMuI7QAmS0hMY5gi3JPdOgMqUbl%2BXwIIkIvB%2Bya16EVLqjLE3JZDrM6Hb%2FhY1diYOliU%2B9rUhf5Rqn4iNIzSy4VHi%2FngVyX2TRBVMOGJzo6cHKQIi38oOSZTZZmF5v1SL4rg9ERSgyXT0m2%2FreTHyXAOTmgLgZ%2B4npS6ueOJlxx0yMw3pbPzHEBR%2B%2Fgjy3TSihYb6hrqaPVWkhEgwzeC0vwY6pQFdTsf4utw9l4MrXTMu06gv77blv2asxJTAHi42HO4Aa4r4zWWUp3RkNdzR5oiYxruoJWMPREBzGPJmVGqddRMQnG2X9LCdjJa%2B%2F7FUIYYy9RRfValKRrnzlnAl5deBK8jMu%2BmVTAffcCGPoZJ152jBY1FqgIQIlox3nEEJAlkoZ33NURHXpTcNAJaZOmU7LQopPljQ0jf9gWnJ6QP1ZHzRAYs6zUs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2000%3A00%3A57%20GMT&X-Amz-Signature=cdeac8be42f85a9e3d3bb47df834808fc4178f06d9d3cc435e2042dde6da33ea'
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEF5NZIN5%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T160057Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG0aCXVzLXdlc3QtMiJHMEUCICQZwBitIu4Neglo3bHNBCT1E8PvFZ%2FQxlTfmkY3HJ4KAiEAg6tKCsnDIROPImOf76nbWN6GO9m4dMJNp1nqHUL9bjsq5gMI1v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDHLIT5JpxqhsIaJq6iq6AzJNAysmsomDmw%2B4G9FYLibD3eb%2Fe704hzNOoDYLc%2B7TRT4zkrhXlK2nwowzVJQdWgN%2BaJw0CXT3sCavmv3mHNKRaGbJ8k8xsskjo05Ge%2FGEc88ZMGDiGcWroBAlx3QvC7B0uSFXY%2BnyC28u6VSxj3PyAkhPd45oMhiaDJqi6b3zmq19EpJfVGdW7F7ihIPserIIo6g9at3uD3uVKeAxlKh%2B0aNHIEeb3QsUpzjUJLT0k3Yr3RzlFwWC%2B8a4kmAFK12kmir3%2FnkP9lPSD8%2FJ3zeT3xqrzaQ%2BhJwxI9IcJ88lWHMpHTNTipcI2XuPlee9kENxebUJehDSPNGKYmho2UxmXF79E0aLRDyC6c8FOpoli1LX7%2ByPQjR8%2BphqtIewMuI7QAmS0hMY5gi3JPdOgMqUbl%2BXwIIkIvB%2Bya16EVLqjLE3JZDrM6Hb%2FhY1diYOliU%2B9rUhf5Rqn4iNIzSy4VHi%2FngVyX2TRBVMOGJzo6cHKQIi38oOSZTZZmF5v1SL4rg9ERSgyXT0m2%2FreTHyXAOTmgLgZ%2B4npS6ueOJlxx0yMw3pbPzHEBR%2B%2Fgjy3TSihYb6hrqaPVWkhEgwzeC0vwY6pQFdTsf4utw9l4MrXTMu06gv77blv2asxJTAHi42HO4Aa4r4zWWUp3RkNdzR5oiYxruoJWMPREBzGPJmVGqddRMQnG2X9LCdjJa%2B%2F7FUIYYy9RRfValKRrnzlnAl5deBK8jMu%2BmVTAffcCGPoZJ152jBY1FqgIQIlox3nEEJAlkoZ33NURHXpTcNAJaZOmU7LQopPljQ0jf9gWnJ6QP1ZHzRAYs6zUs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2000%3A00%3A57%20GMT&X-Amz-Signature=f1ef4ff26fe06993f38cd788ec1fef1567a7354ddc4a0c8139d17f39de218934'
$ws.Range("L3").Style = "Normal"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = '281474991205262-1743535543888'
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = 'Mobile Usage'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '2025-04-01T13:25:43.888'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '281474991205262'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '132'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '52215867'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = 'EMMANUEL SALCEDO'
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 20.64971861
$ws.Range("I4").Value = -103.308130149
$ws.Range("J4").Value = 0
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEF5NZIN5%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T160057Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG0aCXVzLXdlc3QtMiJHMEUCICQZwBitIu4Neglo3bHNBCT1E8PvFZ%2FQxlTfmkY3HJ4KAiEAg6tKCsnDIROPImOf76nbWN6GO9m4dMJNp1nqHUL9bjsq5gMI1v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDHLIT5JpxqhsIaJq6iq6AzJNAysmsomDmw%2B4G9FYLibD3eb%2Fe704hzNOoDYLc%2B7TRT4zkrhXlK2nwowzVJQdWgN%2BaJw0CXT3sCavmv3mHNKRaGbJ8k8xsskjo05Ge%2FGEc88ZMGDiGcWroBAlx3QvC7B0uSFXY%2BnyC28u6VSxj3PyAkhPd45oMhiaDJqi6b3zmq19EpJfVGdW7F7ihIPserIIo6g9at3uD3uVKeAxlKh%2B0aNHIEeb3QsUpzjUJLT0k3Yr3RzlFwWC%2B8a4kmAFK12kmir3%2FnkP9lPSD8%2FJ3zeT3xqrzaQ%2BhJwxI9IcJ88lWHMpHTNTipcI2XuPlee9kENxebUJehDSPNGKYmho2UxmXF79E0aLRDyC6c8FOpoli1LX7%2ByPQjR8%2BphqtIewMuI7QAmS0hMY5gi3JPdOgMqUbl%2BXwIIkIvB%2Bya16EVLqjLE3JZDrM6Hb%2FhY1diYOliU%2B9rUhf5Rqn4iNIzSy4VHi%2FngVyX2TRBVMOGJzo6cHKQIi38oOSZTZZmF5v1SL4rg9ERSgyXT0m2%2FreTHyXAOTmgLgZ%2B4npS6ueOJlxx0yMw3pbPzHEBR%2B%2Fgjy3TSihYb6hrqaPVWkhEgwzeC0vwY6pQFdTsf4utw9l4MrXTMu06gv77blv2asxJTAHi42HO4Aa4r4zWWUp3RkNdzR5oiYxruoJWMPREBzGPJmVGqddRMQnG2X9LCdjJa%2B%2F7FUIYYy9RRfValKRrnzlnAl5deBK8jMu%2BmVTAffcCGPoZJ152jBY1FqgIQIlox3nEEJAlkoZ33NURHXpTcNAJaZOmU7LQopPljQ0jf9gWnJ6QP1ZHzRAYs6zUs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2000%3A00%3A57%20GMT&X-Amz-Signature=df70eb1c7d1a2fecf3f8f995f44e181f954e62ffedd3ba6f94c4b604a7715c07'
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = 'No video URL'
$ws.Range("L4").Style = "Normal"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = '281474991205262-1743532928226'
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'Mobile Usage'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '2025-04-01T12:42:08.226'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '281474991205262'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '132'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '52215867'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = 'EMMANUEL SALCEDO'
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 20.65968375
$ws.Range("I5").Value = -103.29894644
$ws.Range("J5").Value = 0
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEF5NZIN5%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T160057Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG0aCXVzLXdlc3QtMiJHMEUCICQZwBitIu4Neglo3bHNBCT1E8PvFZ%2FQxlTfmkY3HJ4KAiEAg6tKCsnDIROPImOf76nbWN6GO9m4dMJNp1nqHUL9bjsq5gMI1v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDHLIT5JpxqhsIaJq6iq6AzJNAysmsomDmw%2B4G9FYLibD3eb%2Fe704hzNOoDYLc%2B7TRT4zkrhXlK2nwowzVJQdWgN%2BaJw0CXT3sCavmv3mHNKRaGbJ8k8xsskjo05Ge%2FGEc88ZMGDiGcWroBAlx3QvC7B0uSFXY%2BnyC28u6VSxj3PyAkhPd45oMhiaDJqi6b3zmq19EpJfVGdW7F7ihIPserIIo6g9at3uD3uVKeAxlKh%2B0aNHIEeb3QsUpzjUJLT0k3Yr3RzlFwWC%2B8a4kmAFK12kmir3%2FnkP9lPSD8%2FJ3zeT3xqrzaQ%2BhJwxI9IcJ88lWHMpHTNTipcI2XuPlee9kENxebUJehDSPNGKYmho2UxmXF79E0aLRDyC6c8FOpoli1LX7%2ByPQjR8%2BphqtIewMuI7QAmS0hMY5gi3JPdOgMqUbl%2BXwIIkIvB%2Bya16EVLqjLE3JZDrM6Hb%2FhY1diYOliU%2B9rUhf5Rqn4iNIzSy4VHi%2FngVyX2TRBVMOGJzo6cHKQIi38oOSZTZZmF5v1SL4rg9ERSgyXT0m2%2FreTHyXAOTmgLgZ%2B4npS6ueOJlxx0yMw3pbPzHEBR%2B%2Fgjy3TSihYb6hrqaPVWkhEgwzeC0vwY6pQFdTsf4utw9l4MrXTMu06gv77blv2asxJTAHi42HO4Aa4r4zWWUp3RkNdzR5oiYxruoJWMPREBzGPJmVGqddRMQnG2X9LCdjJa%2B%2F7FUIYYy9RRfValKRrnzlnAl5deBK8jMu%2BmVTAffcCGPoZJ152jBY1FqgIQIlox3nEEJAlkoZ33NURHXpTcNAJaZOmU7LQopPljQ0jf9gWnJ6QP1ZHzRAYs6zUs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2000%3A00%3A57%20GMT&X-Amz-Signature=845576b66372b83bc80af75e2da037208b2e477f40aefb0be42868e42214acc8'
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = 'No video URL'
$ws.Range("L5").Style = "Normal"

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = '281474991205262-1743532888631'
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'Mobile Usage'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '2025-04-01T12:41:28.631'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '281474991205262'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '132'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '52215867'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = 'EMMANUEL SALCEDO'
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 20.661818499
$ws.Range("I6").Value = -103.298643339
$ws.Range("J6").Value = 0
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEF5NZIN5%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T160057Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG0aCXVzLXdlc3QtMiJHMEUCICQZwBitIu4Neglo3bHNBCT1E8PvFZ%2FQxlTfmkY3HJ4KAiEAg6tKCsnDIROPImOf76nbWN6GO9m4dMJNp1nqHUL9bjsq5gMI1v%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDHLIT5JpxqhsIaJq6iq6AzJNAysmsomDmw%2B4G9FYLibD3eb%2Fe704hzNOoDYLc%2B7TRT4zkrhXlK2nwowzVJQdWgN%2BaJw0CXT3sCavmv3mHNKRaGbJ8k8xsskjo05Ge%2FGEc88ZMGDiGcWroBAlx3QvC7B0uSFXY%2BnyC28u6VSxj3PyAkhPd45oMhiaDJqi6b3zmq19EpJfVGdW7F7ihIPserIIo6g9at3uD3uVKeAxlKh%2B0aNHIEeb3QsUpzjUJLT0k3Yr3RzlFwWC%2B8a4kmAFK12kmir3%2FnkP9lPSD8%2FJ3zeT3xqrzaQ%2BhJwxI9IcJ88lWHMpHTNTipcI2XuPlee9kENxebUJehDSPNGKYmho2UxmXF79E0aLRDyC6c8FOpoli1LX7%2ByPQjR8%2BphqtIewMuI7QAmS0hMY5gi3JPdOgMqUbl%2BXwIIkIvB%2Bya16EVLqjLE3JZDrM6Hb%2FhY1diYOliU%2B9rUhf5Rqn4iNIzSy4VHi%2FngVyX2TRBVMOGJzo6cHKQIi38oOSZTZZmF5v1SL4rg9ERSgyXT0m2%2FreTHyXAOTmgLgZ%2B4npS6ueOJlxx0yMw3pbPzHEBR%2B%2Fgjy3TSihYb6hrqaPVWkhEgwzeC0vwY6pQFdTsf4utw9l4MrXTMu06gv77blv2asxJTAHi42HO4Aa4r4zWWUp3RkNdzR5oiYxruoJWMPREBzGPJmVGqddRMQnG2X9LCdjJa%2B%2F7FUIYYy9RRfValKRrnzlnAl5deBK8jMu%2BmVTAffcCGPoZJ152jBY1FqgIQIlox3nEEJAlkoZ33NURHXpTcNAJaZOmU7LQopPljQ0jf9gWnJ6QP1ZHzRAYs6zUs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2000%3A00%3A57%20GMT&X-Amz-Signature=af56c8e3bbe20a65777dbb2d9a37885deaab0c0aeb0125a0213b7741e068442b'
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = 'No video URL'
$ws.Range("L6").Style = "Normal"
